$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "50.990.70"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.929.87"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.28%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "377.38"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.93"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.42%  "
$ws.Range("E7").Value = "  -1.65%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.583"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.48"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.63%  "
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0837"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.399.46"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.97"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.35"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.914.25"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.975"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "50.945.43"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.16"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -10.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.09"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.50"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -5.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0950"
$ws.Range("D22").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.25"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.38"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.89"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.94%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.19"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +8.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.62"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.28%  "
$ws.Range("E28").Value = "  +9.52%  "
$ws.Range("E29").Value = "  -2.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.54"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.77"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "33.96"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "50.57"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0452"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.17%  "
$ws.Range("E36").Value = "  -1.80%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.96"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.56%  "
$ws.Range("E39").Value = "  -3.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.53"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -5.87%  "
$ws.Range("E41").Value = "  -1.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.77"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.53"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.19"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -6.24%  "
$ws.Range("E45").Value = "  -2.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.273"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.21%  "
$ws.Range("E47").Value = "  -3.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.002.45"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.21"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.08%  "
$ws.Range("E50").Value = "  -1.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.481"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +12.28%  "
